$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the header style (bold,
# centered, bordered) used by the other header cells (e.g. E1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row.
$timestamps = @(
    "2021-10-05 13:40:11.079251",
    "2021-10-05 13:40:11.079264",
    "2021-10-05 13:40:11.079268",
    "2021-10-05 13:40:11.079271",
    "2021-10-05 13:40:11.079275",
    "2021-10-05 13:40:11.079278",
    "2021-10-05 13:40:11.079281",
    "2021-10-05 13:40:11.079284",
    "2021-10-05 13:40:11.079287",
    "2021-10-05 13:40:11.079290",
    "2021-10-05 13:40:11.079293",
    "2021-10-05 13:40:11.079297",
    "2021-10-05 13:40:11.079300",
    "2021-10-05 13:40:11.079303",
    "2021-10-05 13:40:11.079306"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
